$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J45").Value = 0.2388379152847414
$ws.Range("I46").Value = 0.3744780054549828
$ws.Range("H47").Value = 0.1336718235993181
$ws.Range("G48").Value = 0.08834060834722172
$ws.Range("F49").Value = 0.02147918641116785
$ws.Range("E50").Value = -0.00810701594554874
$ws.Range("D51").Value = -0.02625767267518964
$ws.Range("C52").Value = -0.04428949692388896
$ws.Range("B53").Value = -0.09587373626955231
